{"js": "// Replace each two-digit multiplication prompt in the practice table with\n// its new value. Every occurrence is a unique literal string (the \"a\u00d7b=\"\n// prompts), so a simple search + replace-in-place keeps all paragraph /\n// run formatting (font, size, alignment) untouched \u2014 only the <w:t> text\n// itself changes, matching the source diff exactly.\nconst replacements = [\n  [\"27\u00d749=\", \"61\u00d742=\"],\n  [\"25\u00d763=\", \"28\u00d756=\"],\n  [\"44\u00d776=\", \"65\u00d794=\"],\n  [\"55\u00d784=\", \"24\u00d711=\"],\n  [\"32\u00d775=\", \"47\u00d731=\"],\n  [\"61\u00d727=\", \"91\u00d761=\"],\n  [\"75\u00d726=\", \"41\u00d761=\"],\n  [\"60\u00d731=\", \"72\u00d756=\"],\n  [\"31\u00d795=\", \"66\u00d761=\"],\n  [\"27\u00d761=\", \"65\u00d711=\"],\n  [\"91\u00d740=\", \"94\u00d764=\"],\n  [\"41\u00d768=\", \"29\u00d760=\"],\n  [\"12\u00d735=\", \"28\u00d747=\"],\n  [\"46\u00d764=\", \"31\u00d728=\"],\n  [\"54\u00d731=\", \"46\u00d717=\"],\n  [\"55\u00d712=\", \"65\u00d719=\"],\n  [\"60\u00d775=\", \"59\u00d727=\"],\n  [\"58\u00d753=\", \"33\u00d772=\"],\n  [\"19\u00d766=\", \"78\u00d752=\"],\n  [\"36\u00d715=\", \"68\u00d796=\"],\n  [\"63\u00d756=\", \"53\u00d797=\"],\n  [\"28\u00d728=\", \"68\u00d762=\"],\n  [\"49\u00d791=\", \"78\u00d766=\"],\n  [\"50\u00d792=\", \"42\u00d784=\"],\n  [\"51\u00d774=\", \"68\u00d721=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit multiplication prompt in the practice table with\n# its new value. Every prompt text is unique within the document, so a\n# whole-document Find/Replace (wdReplaceOne) per pair touches only the\n# matching run's text and leaves paragraph/run formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"27\u00d749=\", \"61\u00d742=\"),\n    @(\"25\u00d763=\", \"28\u00d756=\"),\n    @(\"44\u00d776=\", \"65\u00d794=\"),\n    @(\"55\u00d784=\", \"24\u00d711=\"),\n    @(\"32\u00d775=\", \"47\u00d731=\"),\n    @(\"61\u00d727=\", \"91\u00d761=\"),\n    @(\"75\u00d726=\", \"41\u00d761=\"),\n    @(\"60\u00d731=\", \"72\u00d756=\"),\n    @(\"31\u00d795=\", \"66\u00d761=\"),\n    @(\"27\u00d761=\", \"65\u00d711=\"),\n    @(\"91\u00d740=\", \"94\u00d764=\"),\n    @(\"41\u00d768=\", \"29\u00d760=\"),\n    @(\"12\u00d735=\", \"28\u00d747=\"),\n    @(\"46\u00d764=\", \"31\u00d728=\"),\n    @(\"54\u00d731=\", \"46\u00d717=\"),\n    @(\"55\u00d712=\", \"65\u00d719=\"),\n    @(\"60\u00d775=\", \"59\u00d727=\"),\n    @(\"58\u00d753=\", \"33\u00d772=\"),\n    @(\"19\u00d766=\", \"78\u00d752=\"),\n    @(\"36\u00d715=\", \"68\u00d796=\"),\n    @(\"63\u00d756=\", \"53\u00d797=\"),\n    @(\"28\u00d728=\", \"68\u00d762=\"),\n    @(\"49\u00d791=\", \"78\u00d766=\"),\n    @(\"50\u00d792=\", \"42\u00d784=\"),\n    @(\"51\u00d774=\", \"68\u00d721=\")\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    [void]$rng.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
